$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note in column A, row 3 (the text that was overflowing into the
# nurse-details columns for admins)
$ws.Range("A3").Value = "For sis Manyi's laptop"

# Give column A enough width to show the new text
$ws.Columns.Item(1).ColumnWidth = 27.6328125

# Update the selection to match what was saved in the source workbook
$ws.Range("E3:F3").Select()
